# Rename the embedded logo pictures:
#   - the two Pearson "PearsonLogo.png" inline pictures (footers) go from
#     name "image1.png" -> "image2.png"
#   - the "BTec_Logo-Orange" inline picture (header) goes from
#     name "image2.jpg" -> "image1.jpg"
#
# InlineShape.Name maps to the picture's wp:docPr/@name (and, where the
# runtime supports it, the mirrored pic:cNvPr/@name) so we drive the
# rename entirely through the Word object model - no raw XML poking.

$d = $word.ActiveDocument
$sec = $d.Sections(1)

function Rename-FirstInlineShape($rng, $newName) {
    $shp = $rng.InlineShapes(1)
    $shp.Name = $newName
}

# --- Footers: both PearsonLogo pictures, image1.png -> image2.png -------
for ($fi = 1; $fi -le $sec.Footers.Count; $fi++) {
    $f = $sec.Footers($fi)
    if (-not $f.Exists) { continue }

    if ($f.Range.InlineShapes.Count -gt 0) {
        Rename-FirstInlineShape $f.Range "image2.png"
        continue
    }

    # Some footers only resolve their picture through a specific
    # paragraph's range rather than the footer's full range.
    $paraCount = $f.Range.Paragraphs.Count
    for ($pi = 1; $pi -le $paraCount; $pi++) {
        $p = $f.Range.Paragraphs($pi)
        if ($p.Range.InlineShapes.Count -gt 0) {
            Rename-FirstInlineShape $p.Range "image2.png"
        }
    }
}

# --- Header: BTec_Logo-Orange picture, image2.jpg -> image1.jpg ---------
for ($hi = 1; $hi -le $sec.Headers.Count; $hi++) {
    $h = $sec.Headers($hi)
    if (-not $h.Exists) { continue }

    if ($h.Range.InlineShapes.Count -gt 0) {
        Rename-FirstInlineShape $h.Range "image1.jpg"
        continue
    }

    $paraCount = $h.Range.Paragraphs.Count
    for ($pi = 1; $pi -le $paraCount; $pi++) {
        $p = $h.Range.Paragraphs($pi)
        if ($p.Range.InlineShapes.Count -gt 0) {
            Rename-FirstInlineShape $p.Range "image1.jpg"
        }
    }
}
